# Add a PriceChange/UpDown result to row 6 (which was previously the last
# row and was missing its X/Y values), then append a brand new row 7 of
# scan data - matching a fresh day's worth of repeater-scanner output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-missing X6 / Y6 cells -------------------------
$ws.Range("X6").Value = 0.059999000000001246
$ws.Range("Y6").Value = "Up"

# --- Append new row 7 ------------------------------------------------------
$ws.Range("A7").Value = 42648.886597222219
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"

$ws.Range("B7").Value = 11
$ws.Range("C7").Value = "Buy"
$ws.Range("D7").Value = 52
$ws.Range("E7").Value = 10935
$ws.Range("F7").Value = 631
$ws.Range("G7").Value = 66
$ws.Range("H7").Value = 31
$ws.Range("I7").Value = 92
$ws.Range("J7").Value = 7
$ws.Range("K7").Value = 20010
$ws.Range("L7").Value = 158
$ws.Range("M7").Value = 76
$ws.Range("N7").Value = 49
$ws.Range("O7").Value = 4
$ws.Range("P7").Value = "Named"
$ws.Range("Q7").Value = 41.162214763508182
$ws.Range("R7").Value = 0

$ws.Range("S7").Value = 0.0616
$ws.Range("S7").NumberFormat = "0.00%"
$ws.Range("T7").Value = -0.032
$ws.Range("T7").NumberFormat = "0.00%"

$ws.Range("U7").Value = 2.2599999999999998
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = 0
